$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("table3_COPR_s_g_hs_inc_raw")
$ws1.Range("D5").Value = 0
$ws1.Range("B6").Value = 5
$ws1.Range("D6").Value = 5
$ws1.Range("B7").Value = 16
$ws1.Range("D7").Value = 16

$ws2 = $wb.Worksheets.Item("table3_COPR_s_g_hs_inc_divtot")
$ws2.Range("B2").Value = 0.25
$ws2.Range("D2").Value = 0.56200000000000006
$ws2.Range("B3").Value = 0.312
$ws2.Range("D4").Value = 0.125
$ws2.Range("B5").Value = 0.125
$ws2.Range("D5").Value = 0
$ws2.Range("B6").Value = 0.312
$ws2.Range("D6").Value = 0.312
$ws2.Range("B7").Value = 16
$ws2.Range("D7").Value = 16

$ws3 = $wb.Worksheets.Item("table3_COPR_s_g_hs_inc_divext")
$ws3.Range("D2").Value = 0.81799999999999995
$ws3.Range("D4").Value = 0.182
$ws3.Range("D5").Value = 0
$ws3.Range("B6").Value = 5
$ws3.Range("D6").Value = 5
$ws3.Range("B7").Value = 16
$ws3.Range("D7").Value = 16

$ws4 = $wb.Worksheets.Item("table3_COPR_s_g_hs_dec_raw")
$ws4.Range("D5").Value = 11
$ws4.Range("B6").Value = 5
$ws4.Range("D6").Value = 5
$ws4.Range("B7").Value = 16
$ws4.Range("D7").Value = 16

$ws5 = $wb.Worksheets.Item("table3_COPR_s_g_hs_dec_divtot")
$ws5.Range("B3").Value = 0.062
$ws5.Range("B5").Value = 0.625
$ws5.Range("D5").Value = 0.68799999999999994
$ws5.Range("B6").Value = 0.312
$ws5.Range("D6").Value = 0.312
$ws5.Range("B7").Value = 16
$ws5.Range("D7").Value = 16

$ws6 = $wb.Worksheets.Item("table3_COPR_s_g_hs_dec_divext")
$ws6.Range("B6").Value = 5
$ws6.Range("D6").Value = 5
$ws6.Range("B7").Value = 16
$ws6.Range("D7").Value = 16

$ws1.Activate()
